# Update "想去人数" (F column) figures and mark G2 as "不可售" on both the
# "展览" sheet and the "全部类型" sheet (they carry duplicate data).

$wb = $excel.ActiveWorkbook

# Map of F-column row -> new numeric value to apply on both sheets.
$fUpdates = @{
    3  = 692
    4  = 690
    5  = 245
    6  = 16
    7  = 1088
    9  = 1615
    10 = 5713
    11 = 467
    12 = 317
    13 = 258
    14 = 70
    15 = 295
    17 = 4649
    18 = 240
    19 = 1229
    20 = 126
    21 = 93
    23 = 88
    24 = 228
    25 = 80
    28 = 362
    29 = 47
    31 = 72
    32 = 25
    33 = 39
    34 = 11
    35 = 51
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $fUpdates.Keys) {
        $ws.Range("F$row").Value = $fUpdates[$row]
    }

    # G2: 65 -> "不可售"
    $ws.Range("G2").Value = "不可售"
}
